$wb = $excel.ActiveWorkbook

# --- "Test Data" sheet: update Id and CreatedDate for row 2 ---
$wsTest = $wb.Worksheets.Item("Test Data")
$wsTest.Range("A2").Value = 56
$wsTest.Range("G2").Value = 44041.472538229165

# --- "Measurement Data" sheet: update Id, measured Current*, CreatedDate, TestId per row ---
$wsMeas = $wb.Worksheets.Item("Measurement Data")
# Row 2
$wsMeas.Range("A2").Value = 346
$wsMeas.Range("T2").Value = 3.39566
$wsMeas.Range("Y2").Value = 44041.47266554398
$wsMeas.Range("Z2").Value = 56
# Row 3
$wsMeas.Range("A3").Value = 347
$wsMeas.Range("T3").Value = 1.2619
$wsMeas.Range("Y3").Value = 44041.47277565972
$wsMeas.Range("Z3").Value = 56
# Row 4
$wsMeas.Range("A4").Value = 348
$wsMeas.Range("T4").Value = 0.05568
$wsMeas.Range("Y4").Value = 44041.47287704861
$wsMeas.Range("Z4").Value = 56
# Row 5
$wsMeas.Range("A5").Value = 349
$wsMeas.Range("T5").Value = 0.01723
$wsMeas.Range("Y5").Value = 44041.472987650464
$wsMeas.Range("Z5").Value = 56
# Row 6
$wsMeas.Range("A6").Value = 350
$wsMeas.Range("Y6").Value = 44041.4730900463
$wsMeas.Range("Z6").Value = 56
# Row 7
$wsMeas.Range("A7").Value = 351
$wsMeas.Range("U7").Value = 2.98841
$wsMeas.Range("Y7").Value = 44041.47321643519
$wsMeas.Range("Z7").Value = 56
# Row 8
$wsMeas.Range("A8").Value = 352
$wsMeas.Range("U8").Value = 1.50098
$wsMeas.Range("Y8").Value = 44041.47332596065
$wsMeas.Range("Z8").Value = 56
# Row 9
$wsMeas.Range("A9").Value = 353
$wsMeas.Range("U9").Value = 0.30014
$wsMeas.Range("Y9").Value = 44041.47342638889
$wsMeas.Range("Z9").Value = 56
# Row 10
$wsMeas.Range("A10").Value = 354
$wsMeas.Range("U10").Value = 0.08002
$wsMeas.Range("Y10").Value = 44041.47353634259
$wsMeas.Range("Z10").Value = 56
# Row 11
$wsMeas.Range("A11").Value = 355
$wsMeas.Range("Y11").Value = 44041.473637696756
$wsMeas.Range("Z11").Value = 56
# Row 12
$wsMeas.Range("A12").Value = 356
$wsMeas.Range("V12").Value = 2.97702
$wsMeas.Range("Y12").Value = 44041.47375428241
$wsMeas.Range("Z12").Value = 56
# Row 13
$wsMeas.Range("A13").Value = 357
$wsMeas.Range("V13").Value = 1.17892
$wsMeas.Range("Y13").Value = 44041.47386388889
$wsMeas.Range("Z13").Value = 56
# Row 14
$wsMeas.Range("A14").Value = 358
$wsMeas.Range("V14").Value = 0.09261
$wsMeas.Range("Y14").Value = 44041.47396458333
$wsMeas.Range("Z14").Value = 56
# Row 15
$wsMeas.Range("A15").Value = 359
$wsMeas.Range("V15").Value = 0.02497
$wsMeas.Range("Y15").Value = 44041.47406550926
$wsMeas.Range("Z15").Value = 56
# Row 16
$wsMeas.Range("A16").Value = 360
$wsMeas.Range("Y16").Value = 44041.47417581019
$wsMeas.Range("Z16").Value = 56
# Row 17
$wsMeas.Range("A17").Value = 361
$wsMeas.Range("W17").Value = 2.87782
$wsMeas.Range("Y17").Value = 44041.47430158565
$wsMeas.Range("Z17").Value = 56
# Row 18
$wsMeas.Range("A18").Value = 362
$wsMeas.Range("W18").Value = 1.24021
$wsMeas.Range("Y18").Value = 44041.47441096065
$wsMeas.Range("Z18").Value = 56
# Row 19
$wsMeas.Range("A19").Value = 363
$wsMeas.Range("W19").Value = 0.10708
$wsMeas.Range("Y19").Value = 44041.47451184028
$wsMeas.Range("Z19").Value = 56
# Row 20
$wsMeas.Range("A20").Value = 364
$wsMeas.Range("W20").Value = 0.0299
$wsMeas.Range("Y20").Value = 44041.474612997685
$wsMeas.Range("Z20").Value = 56
# Row 21
$wsMeas.Range("A21").Value = 365
$wsMeas.Range("Y21").Value = 44041.47481605324
$wsMeas.Range("Z21").Value = 56
